$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45202 -> 2023-10-03).
# Every row from 2 to 171 needs to be bumped by one day to 45203 (2023-10-04).
$ws.Range("C2:C171").Value = 45203
